# Update the LinkedIn carousel deck: the news item changes from the
# "NTPC Green Energy commissions Bhadla Solar Project" story to the
# "NTPC Green Energy and GAIL joint venture" story. Every slide's title
# (shape 1) is refreshed to the new headline, and each slide's content
# placeholder (shape 2) has its two detail bullets (paragraphs 2 and 3 -
# paragraph 1 is the "image missing" placeholder line) replaced with new
# supporting points.

$p = $ppt.ActivePresentation

$newTitle = "NTPC Green Energy And GAIL Join Hands To Launch 50:50 Joint Venture For Renewable Energy Projects - SolarQuarter"

$bullets = @{
    1 = @(
        "NTPC Green Energy and GAIL have formed a 50:50 joint venture.",
        "The joint venture focuses on renewable energy projects."
    )
    2 = @(
        "The collaboration aims to enhance the renewable energy capacity in India.",
        "Both companies are leaders in their respective sectors."
    )
    3 = @(
        "The joint venture will leverage NTPC's expertise in power generation.",
        "GAIL will contribute its experience in energy infrastructure."
    )
    4 = @(
        "The initiative aligns with India's commitment to increasing renewable energy sources.",
        "The joint venture is expected to support government policies on clean energy."
    )
    5 = @(
        "The partnership is part of a broader strategy to reduce carbon emissions.",
        "It aims to contribute to sustainable development goals."
    )
    6 = @(
        "The joint venture will explore various renewable energy technologies.",
        "It signifies a strategic alliance in the energy sector."
    )
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    # Shape 1: title
    $titleShape = $slide.Shapes.Item(1)
    $titleShape.TextFrame.TextRange.Runs(1).Text = $newTitle

    # Shape 2: content placeholder - paragraph 1 is the image-missing
    # placeholder, paragraphs 2 and 3 are the two detail bullets.
    $bodyShape = $slide.Shapes.Item(2)
    $bodyTr = $bodyShape.TextFrame.TextRange
    $points = $bullets[$i]
    $bodyTr.Paragraphs(2).Runs(1).Text = $points[0]
    $bodyTr.Paragraphs(3).Runs(1).Text = $points[1]
}
